$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harvester (column B) and experimentDesign (column D) for all data rows (2-13)
for ($r = 2; $r -le 13; $r++) {
    $ws.Range("B$r").Value = "S.GISH"
    $ws.Range("D$r").Value = "90minuteInduction"
}

# Update strain (column F) per block of rows
for ($r = 2; $r -le 4; $r++) {
    $ws.Range("F$r").Value = "KN99alpha"
}
for ($r = 8; $r -le 10; $r++) {
    $ws.Range("F$r").Value = "TDY1700"
}
for ($r = 11; $r -le 13; $r++) {
    $ws.Range("F$r").Value = "TDY1319"
}

# Restore last-used selection as recorded in the saved workbook
$ws.Range("F12:F13").Select()
